# Apply "Despues de las pruebas de Agaela" edits across all board sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: inicial ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Value = "TAB. RÁPIDO"
$ws.Range("B1").Value = "TAB. VERBOS"
$ws.Range("C1").Value = "TAB. COMIDA"
$ws.Range("D1").Value = "TAB. OBJETOS"
$ws.Range("A2").Value = "TAB. PERSONAS"
$ws.Range("B2").Value = "TAB. LUGARES"
$ws.Range("C2").Value = "TAB. TRANSPORTE"
$ws.Range("D2").Value = "TAB. CASA"
$ws.Range("A3").Value = "TAB. ANIMALES"
$ws.Range("B3").Value = "TAB. CUERPO"
$ws.Range("C3").Value = "TAB. CONCEPTOS"
$ws.Range("D3").Value = "?"

# --- Sheet 2: rápido ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A1").Value = "TAB. INICIAL"
$ws.Range("B1").Value = "QUERER"
$ws.Range("C1").Value = "NECESITAR"
$ws.Range("D1").Value = "AYUDA"
$ws.Range("A2").Value = "SÍ"
$ws.Range("B2").Value = "NO"
$ws.Range("C2").Value = "PARA"
$ws.Range("D2").Value = "COMER"
$ws.Range("A3").Value = "YO"
$ws.Range("B3").Value = "TÚ"
$ws.Range("C3").Value = "BIEN"
$ws.Range("D3").Value = "MAL"
$ws.Rows.Item(4).Delete()
$ws.Range("B5").Select()

# --- Sheet 3: verbos ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A1").Value = "TAB. INICIAL"
$ws.Range("B1").Value = "TAB. VERBOS2"
$ws.Range("C1").Value = "VER"
$ws.Range("D1").Value = "IR"
$ws.Range("A2").Value = "ABRIR"
$ws.Range("B2").Value = "VIAJAR"
$ws.Range("C2").Value = "DORMIR"
$ws.Range("D2").Value = "PODER"
$ws.Range("A3").Value = "NECESITAR"
$ws.Range("B3").Value = "QUERER"
$ws.Range("C3").Value = "BEBER"
$ws.Range("D3").Value = "COMER"
$ws.Rows.Item(4).Delete()
$ws.Range("B2").Select()

# --- Sheet 4: verbos2 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("A1").Value = "TAB. INICIAL"
$ws.Range("B1").Value = "TAB. VERBOS"
$ws.Range("C1").Value = "JUGAR"
$ws.Range("D1").Value = "CERRAR"
$ws.Range("A2").Value = "VOLVER"
$ws.Range("B2").Value = "PENSAR"
$ws.Range("C2").Value = "TENER"
$ws.Range("D2").Value = "DAR"
$ws.Range("A3").Value = "SENTIR"
$ws.Range("B3").Value = "SER"
$ws.Range("C3").Value = "ESCUCHAR"
$ws.Range("D3").Value = "LLEVAR"
$ws.Rows.Item(4).Delete()
$ws.Range("D1").Select()

# --- Sheet 5: comida ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("A1").Value = "TAB. INICIAL"
$ws.Range("B1").Value = "TAB. COMIDA2"
$ws.Range("C1").Value = "LECHUGA"
$ws.Range("D1").Value = "PATATA"
$ws.Range("A2").Value = "PURÉ"
$ws.Range("B2").Value = "SOPA"
$ws.Range("C2").Value = "HABAS"
$ws.Range("D2").Value = "LENTEJAS"
$ws.Range("A3").Value = "TOMATE"
$ws.Range("B3").Value = "CARNE"
$ws.Range("C3").Value = "PESCADO"
$ws.Range("D3").Value = "POLLO"
$ws.Rows.Item(4).Delete()
$ws.Range("A4:D4").Select()

# --- Sheet 6: comida2 ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("A1").Value = "TAB. INICIAL"
$ws.Range("B1").Value = "TAB. COMIDA"
$ws.Range("C1").Value = "PIZZA"
$ws.Range("D1").Value = "HAMBURGUESA"
$ws.Range("A2").Value = "PAN"
$ws.Range("B2").Value = "SANDWICH"
$ws.Range("C2").Value = "QUESO"
$ws.Range("D2").Value = "JAMÓN"
$ws.Range("A3").Value = "CEREALES"
$ws.Range("B3").Value = "GALLETAS"
$ws.Range("C3").Value = "MANZANA"
$ws.Range("D3").Value = "FRESA"
$ws.Rows.Item(4).Delete()
$ws.Range("A4:D4").Select()

# --- Sheet 7: objetos ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("A1").Value = "TAB. INICIAL"
$ws.Range("B1").Value = "DINERO"
$ws.Range("C1").Value = "BATERÍA"
$ws.Range("D1").Value = "SILLA"
$ws.Range("A2").Value = "PARAGUAS"
$ws.Range("B2").Value = "GAFAS"
$ws.Range("C2").Value = "RELOJ"
$ws.Range("D2").Value = "BOLSA"
$ws.Range("A3").Value = "TELEVISIÓN"
$ws.Range("B3").Value = "ORDENADOR"
$ws.Range("C3").Value = "BOTELLA"
$ws.Range("D3").Value = "HORNO"
$ws.Rows.Item(4).Delete()
$ws.Range("A4:D4").Select()

# --- Sheet 8: personas ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("A1").Value = "TAB. INICIAL"
$ws.Range("B1").Value = "YO"
$ws.Range("C1").Value = "NOSOTROS-AS"
$ws.Range("D1").Value = "EL-ELLA"
$ws.Range("A2").Value = "DOCTOR"
$ws.Range("B2").Value = "TÚ"
$ws.Range("C2").Value = "VOSOTROS-AS"
$ws.Range("D2").Value = "ELLOS-ELLAS"
$ws.Range("A3").Value = "PROFESOR"
$ws.Range("B3").Value = "PADRE"
$ws.Range("C3").Value = "MADRE"
$ws.Range("D3").Value = "AMIGO-A"
$ws.Rows.Item(4).Delete()
$ws.Range("D24").Select()

# --- Sheet 9: lugares ---
$ws = $wb.Worksheets.Item(9)
$ws.Range("A1").Value = "TAB. INICIAL"
$ws.Range("B1").Value = "PARQUE"
$ws.Range("C1").Value = "ESCUELA"
$ws.Range("D1").Value = "CASA"
$ws.Range("A2").Value = "PUEBLO"
$ws.Range("B2").Value = "HOSPITAL"
$ws.Range("C2").Value = "CALLE"
$ws.Range("D2").Value = "TIENDA"
$ws.Range("A3").Value = "CINE"
$ws.Range("B3").Value = "TEATRO"
$ws.Range("C3").Value = "RESTAURANTE"
$ws.Range("D3").Value = "MUSEO"
$ws.Rows.Item(4).Delete()
$ws.Range("A5").Select()

# --- Sheet 10: transporte ---
$ws = $wb.Worksheets.Item(10)
$ws.Range("A1").Value = "TAB. INICIAL"
$ws.Range("B1").Value = "AMBULANCIA"
$ws.Range("C1").Value = "CARAVANA"
$ws.Range("D1").Value = "COCHE"
$ws.Range("A2").Value = "AUTOBÚS"
$ws.Range("B2").Value = "AVIÓN"
$ws.Range("C2").Value = "TAXI"
$ws.Range("D2").Value = "METRO"
$ws.Range("A3").Value = "BICICLETA"
$ws.Range("B3").Value = "MOTO"
$ws.Range("C3").Value = "TREN"
$ws.Range("D3").Value = "BARCO"
$ws.Rows.Item(4).Delete()
$ws.Range("B2").Select()

# --- Sheet 11: casa ---
$ws = $wb.Worksheets.Item(11)
$ws.Range("A1").Value = "TAB. INICIAL"
$ws.Range("B1").Value = "BAÑO"
$ws.Range("C1").Value = "COCINA"
$ws.Range("D1").Value = "HABITACIÓN"
$ws.Range("A2").Value = "SALÓN"
$ws.Range("B2").Value = "PASILLO"
$ws.Range("C2").Value = "PUERTA"
$ws.Range("D2").Value = "VENTANA"
$ws.Range("A3").Value = "ASCENSOR"
$ws.Range("B3").Value = "ESCALERA"
$ws.Range("C3").Value = "LLAVE"
$ws.Range("D3").Value = "SOFÁ"
$ws.Rows.Item(4).Delete()
$ws.Range("A4:D4").Select()

# --- Sheet 12: animales ---
$ws = $wb.Worksheets.Item(12)
$ws.Range("A1").Value = "TAB. INICIAL"
$ws.Range("B1").Value = "PERRO"
$ws.Range("C1").Value = "GATO"
$ws.Range("D1").Value = "CABALLO"
$ws.Range("A2").Value = "CERDO"
$ws.Range("B2").Value = "CONEJO"
$ws.Range("C2").Value = "ABEJA"
$ws.Range("D2").Value = "PÁJARO"
$ws.Range("A3").Value = "DELFÍN"
$ws.Range("B3").Value = "TORTUGA"
$ws.Range("C3").Value = "SERPIENTE"
$ws.Range("D3").Value = "MOSQUITO"
$ws.Rows.Item(4).Delete()
$ws.Range("A4:D4").Select()

# --- Sheet 14: conceptos ---
$ws = $wb.Worksheets.Item(14)
$ws.Range("A1").Value = "TAB. INICIAL"
$ws.Range("B1").Value = "IZQUIERDA"
$ws.Range("C1").Value = "DERECHA"
$ws.Range("D1").Value = "MUCHO"
$ws.Range("A2").Value = "POCO"
$ws.Range("B2").Value = "ANTES"
$ws.Range("C2").Value = "DESPUÉS"
$ws.Range("D2").Value = "AHORA"
$ws.Range("A3").Value = "DENTRO"
$ws.Range("B3").Value = "FUERA"
$ws.Range("C3").Value = "CERCA"
$ws.Range("D3").Value = "LEJOS"
$ws.Rows.Item(4).Delete()
$ws.Range("D24").Select()

# --- Sheet 13: cuerpo ---
$ws = $wb.Worksheets.Item(13)
$ws.Range("A1").Value = "TAB. INICIAL"
$ws.Range("B1").Value = "CABEZA"
$ws.Range("C1").Value = "PIE"
$ws.Range("D1").Value = "NARIZ"
$ws.Range("A2").Value = "OREJA"
$ws.Range("B2").Value = "OJO"
$ws.Range("C2").Value = "BOCA"
$ws.Range("D2").Value = "CUELLO"
$ws.Range("A3").Value = "PIERNA"
$ws.Range("B3").Value = "BRAZO"
$ws.Range("C3").Value = "DEDO"
$ws.Range("D3").Value = "MANO"
$ws.Rows.Item(4).Delete()
$ws.Range("C3").Select()
